# Apply the "added information about AU in basic relevane" edit.
$wb = $excel.ActiveWorkbook

# --- 1. Update the keyword-matching question on the "Relevance" sheet (A2) ---
# Add "do know about" to the list of recognized keywords so the new
# "know about AU Andhra University" Q&A pair is reachable.
$wsRelevance = $wb.Worksheets.Item("Relevance")
$wsRelevance.Range("A2").Value = "hi hello hey located established vice chancellor vc registrar rector colleges where au andhra university who how what your smart introduce yourself do know about"

# --- 2. Update the "Basic" sheet ---
$wsBasic = $wb.Worksheets.Item("Basic")

# 2a. Shorten/refresh the existing "established" answer (row 8, column C)
$wsBasic.Range("C8").Value = "Andhra University was established in 1926 by the Madras Act of 1926. The 94-year-old institution is fortunate to have Sir C.R. Reddy as its founder Vice-Chancellor, as the steps taken by this visionary proved to be fruitful in the long run."

# 2b. Add a brand new row (13) with a detailed "know about AU" Q&A entry
$wsBasic.Range("A13").Value = "know about AU Andhra University"
$wsBasic.Range("B13").Value = "Basic"
$wsBasic.Range("C13").Value = "Andhra University is one of the oldest educational institutions in the country and also the first university to be conceived as a residential and teaching-cum-affiliating University, mainly devoted to post-graduate teaching and research. Andhra University was constituted in the year 1926 by the Madras Act of 1926. The 94-year-old institution is fortunate to have Sir C.R. Reddy as its founder Vice-Chancellor, as the steps taken by this visionary proved to be fruitful in the long run. Former President of India Dr. Sarvepalli Radhakrishnan was one of its Vice-Chancellors who succeed Dr. C. R. Reddy in 1931. The leaders of the university always believed that excellence in higher education is the best investment for the country and engaged the services of famous educationists such as Dr. T.R. Seshadri, Dr. S. Bhagavantham, Professor Hiren Mukherjee, Professor Humayan Kabir and Dr. V.K.R.V. Rao, to mention a few who set high standards for teaching and research. Nobel Lariat C V Raman was the proud alumnus of the University and closely associated in laying research foundations in Physics. Padmavibhushan Prof. C R Rao, the renowned statistician of the world, was also the proud alumnus of the University. Ever since its inception in 1926 Andhra University has an impeccable record of catering to the educational needs and solving the sociological problems of the region. The University is relentless in its efforts in maintaining standards in teaching and research, ensuring proper character building and development among the students, encouraging community developmental programmes, nurturing leadership in young men and women and imbibing a sense of responsibility to become good citizens, while striving for excellence in all fronts."

# Match the font colour already used for similar "profile" style text elsewhere
# in the workbook (re-uses the existing cell style rather than creating a new one).
$wsBasic.Range("C13").Font.Color = 2696481

# Keep the sheet's used dimension / page setup tidy like the edited workbook
$wsBasic.PageSetup.Orientation = 1

# Reflect the cell that was being edited/viewed when the change was made
$wsBasic.Range("C8").Select() | Out-Null
